$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.485.48'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Value = '3.482.19'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '591.07'
$ws.Range('E5').Value = '  -1.87%  '
$ws.Range('D6').Value = '178.90'
$ws.Range('E6').Value = '  -1.36%  '
$ws.Range('E7').Value = '  +3.29%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '3.481.01'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E10').Value = '  -2.96%  '
$ws.Range('D11').Value = '6.97'
$ws.Range('E11').Value = '  -2.74%  '
$ws.Range('D12').Value = '0.427'
$ws.Range('E12').Value = '  -2.95%  '
$ws.Range('D13').Value = '4.082.35'
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').Value = '32.16'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').Value = '0.132'
$ws.Range('E15').Value = '  -2.66%  '
$ws.Range('D16').Value = '67.484.93'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('E17').Value = '  -2.80%  '
$ws.Range('D18').Value = '3.480.40'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('E19').Value = '  -3.68%  '
$ws.Range('D20').Value = '14.03'
$ws.Range('E20').Value = '  -3.37%  '
$ws.Range('D21').Value = '382.04'
$ws.Range('E21').Value = '  -4.78%  '
$ws.Range('D22').Value = '7.91'
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('D23').Value = '5.80'
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').Value = '0.536'
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '71.97'
$ws.Range('E26').Value = '  -2.53%  '
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('D28').Value = '10.06'
$ws.Range('E28').Value = '  -3.99%  '
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E31').Value = '  -3.92%  '
$ws.Range('D32').Value = '24.40'
$ws.Range('E32').Value = '  +1.79%  '
$ws.Range('E33').Value = '  -2.86%  '
$ws.Range('E34').Value = '  -4.78%  '
$ws.Range('E35').Value = '  -3.05%  '
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('E37').Value = '  -4.38%  '
$ws.Range('D38').Value = '160.45'
$ws.Range('E38').Value = '  -1.67%  '
$ws.Range('D39').Value = '0.885'
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('D40').Value = '27.97'
$ws.Range('E40').Value = '  +5.90%  '
$ws.Range('E41').Value = '  -3.76%  '
$ws.Range('E42').Value = '  -4.23%  '
$ws.Range('D43').Value = '6.65'
$ws.Range('E43').Value = '  -4.97%  '
$ws.Range('E44').Value = '  -3.80%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.717.49'
$ws.Range('E45').Value = '  -5.98%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '0.0705'
$ws.Range('E46').Value = '  -4.42%  '
$ws.Range('D47').Value = '25.92'
$ws.Range('E47').Value = '  -3.88%  '
$ws.Range('D48').Value = '41.50'
$ws.Range('E48').Value = '  -2.09%  '
$ws.Range('E49').Value = '  -2.74%  '
$ws.Range('D50').Value = '326.84'
$ws.Range('E50').Value = '  -6.36%  '
$ws.Range('E51').Value = '  -3.22%  '
